# Edit: Fri, Jun 12, 2020 12:05:08 AM
#
# 1) Slide 6's table gets a new (built-in) table style applied via the
#    Table Design gallery.
# 2) The presentation's theme colour palette is switched from the
#    "Integral" palette to the standard "Office Theme" palette (the
#    design/theme used by the slide master - and therefore every slide -
#    is changed; fonts and format scheme are left as-is since they were
#    already shared between the two palettes).

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 -------------------------------------------------
$tableSlide = $p.Slides.Item(6)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{EC2AD332-37A3-4239-980F-613CD1900C88}")

# --- 2) Swap the Integral theme colours for the Office Theme colours ----------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = RGB(0, 0, 0)        # dk1
$tcs.Item(2).RGB  = RGB(255, 255, 255)  # lt1
$tcs.Item(3).RGB  = RGB(68, 84, 106)    # dk2
$tcs.Item(4).RGB  = RGB(231, 230, 230)  # lt2
$tcs.Item(5).RGB  = RGB(91, 155, 213)   # accent1
$tcs.Item(6).RGB  = RGB(237, 125, 49)   # accent2
$tcs.Item(7).RGB  = RGB(165, 165, 165)  # accent3
$tcs.Item(8).RGB  = RGB(255, 192, 0)    # accent4
$tcs.Item(9).RGB  = RGB(68, 114, 196)   # accent5
$tcs.Item(10).RGB = RGB(112, 173, 71)   # accent6
$tcs.Item(11).RGB = RGB(5, 99, 193)     # hlink
$tcs.Item(12).RGB = RGB(149, 79, 114)   # folHlink
